$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (501) into the new rows
$srcRow = $ws.Range("A501:V501")
for ($r = 502; $r -le 514; $r++) {
    $dstRow = $ws.Range("A" + $r + ":V" + $r)
    $srcRow.Copy($dstRow)
}

# Write the new training-session data for each player row
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "Entrainement"
$arr[0,1] = 45917
$arr[0,2] = "Global"
$arr[0,3] = "J-3"
$arr[0,4] = "Karahali Souaré"
$arr[0,5] = "right forward"
$arr[0,6] = "01:38:44"
$arr[0,7] = 6.47
$arr[0,8] = 0.75
$arr[0,9] = 5.7
$arr[0,10] = 0.43
$arr[0,11] = 0.22
$arr[0,12] = 0.11
$arr[0,13] = 0.01
$arr[0,14] = 9
$arr[0,15] = 3.59
$arr[0,16] = 31
$arr[0,17] = 5.46
$arr[0,18] = 60
$arr[0,19] = 17
$arr[0,20] = 42
$arr[0,21] = 14
$ws.Range("A502:V502").Value = $arr

$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "Entrainement"
$arr[0,1] = 45917
$arr[0,2] = "Global"
$arr[0,3] = "J-3"
$arr[0,4] = "Levy Ndoutoume"
$arr[0,5] = "left back"
$arr[0,6] = "01:39:57"
$arr[0,7] = 6.63
$arr[0,8] = 0.58
$arr[0,9] = 6.04
$arr[0,10] = 0.42
$arr[0,11] = 0.1
$arr[0,12] = 0.07
$arr[0,13] = 0
$arr[0,14] = 5
$arr[0,15] = 3.56
$arr[0,16] = 29.97
$arr[0,17] = 5.68
$arr[0,18] = 54
$arr[0,19] = 11
$arr[0,20] = 47
$arr[0,21] = 16
$ws.Range("A503:V503").Value = $arr

$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "Entrainement"
$arr[0,1] = 45917
$arr[0,2] = "Global"
$arr[0,3] = "J-3"
$arr[0,4] = "Ilyes Boughanmi"
$arr[0,5] = "center forward"
$arr[0,6] = "01:35:43"
$arr[0,7] = 5.33
$arr[0,8] = 0.46
$arr[0,9] = 4.86
$arr[0,10] = 0.28
$arr[0,11] = 0.16
$arr[0,12] = 0.03
$arr[0,13] = 0
$arr[0,14] = 3
$arr[0,15] = 3.24
$arr[0,16] = 26.83
$arr[0,17] = 5.07
$arr[0,18] = 18
$arr[0,19] = 5
$arr[0,20] = 12
$arr[0,21] = 3
$ws.Range("A504:V504").Value = $arr

$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "Entrainement"
$arr[0,1] = 45917
$arr[0,2] = "Global"
$arr[0,3] = "J-3"
$arr[0,4] = "Omar Benyounes"
$arr[0,5] = "center midfield"
$arr[0,6] = "01:39:21"
$arr[0,7] = 7.18
$arr[0,8] = 0.83
$arr[0,9] = 6.34
$arr[0,10] = 0.61
$arr[0,11] = 0.19
$arr[0,12] = 0.04
$arr[0,13] = 0
$arr[0,14] = 3
$arr[0,15] = 4
$arr[0,16] = 28.68
$arr[0,17] = 4.54
$arr[0,18] = 24
$arr[0,19] = 7
$arr[0,20] = 26
$arr[0,21] = 5
$ws.Range("A505:V505").Value = $arr

$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "Entrainement"
$arr[0,1] = 45917
$arr[0,2] = "Global"
$arr[0,3] = "J-3"
$arr[0,4] = "Jeremie Laurent"
$arr[0,5] = "left forward"
$arr[0,6] = "01:38:54"
$arr[0,7] = 6.42
$arr[0,8] = 0.85
$arr[0,9] = 5.56
$arr[0,10] = 0.4
$arr[0,11] = 0.31
$arr[0,12] = 0.16
$arr[0,13] = 0
$arr[0,14] = 15
$arr[0,15] = 3.85
$arr[0,16] = 29.86
$arr[0,17] = 5.41
$arr[0,18] = 32
$arr[0,19] = 18
$arr[0,20] = 19
$arr[0,21] = 11
$ws.Range("A506:V506").Value = $arr

$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "Entrainement"
$arr[0,1] = 45917
$arr[0,2] = "Global"
$arr[0,3] = "J-3"
$arr[0,4] = "Ilan Ihaddadene"
$arr[0,5] = "center midfield"
$arr[0,6] = "01:34:59"
$arr[0,7] = 7.51
$arr[0,8] = 0.8
$arr[0,9] = 6.7
$arr[0,10] = 0.73
$arr[0,11] = 0.08
$arr[0,12] = 0
$arr[0,13] = 0
$arr[0,14] = 0
$arr[0,15] = 4.66
$arr[0,16] = 24.74
$arr[0,17] = 5.15
$arr[0,18] = 31
$arr[0,19] = 6
$arr[0,20] = 23
$arr[0,21] = 5
$ws.Range("A507:V507").Value = $arr

$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "Entrainement"
$arr[0,1] = 45917
$arr[0,2] = "Global"
$arr[0,3] = "J-3"
$arr[0,4] = "Mattheo Haon"
$arr[0,5] = "right back"
$arr[0,6] = "01:37:41"
$arr[0,7] = 7.74
$arr[0,8] = 0.9
$arr[0,9] = 6.82
$arr[0,10] = 0.55
$arr[0,11] = 0.26
$arr[0,12] = 0.09
$arr[0,13] = 0.01
$arr[0,14] = 7
$arr[0,15] = 4.69
$arr[0,16] = 30.69
$arr[0,17] = 4.46
$arr[0,18] = 35
$arr[0,19] = 7
$arr[0,20] = 24
$arr[0,21] = 9
$ws.Range("A508:V508").Value = $arr

$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "Entrainement"
$arr[0,1] = 45917
$arr[0,2] = "Global"
$arr[0,3] = "J-3"
$arr[0,4] = "Yoan Zouma"
$arr[0,5] = "center back"
$arr[0,6] = "01:37:04"
$arr[0,7] = 5.02
$arr[0,8] = 0.28
$arr[0,9] = 4.73
$arr[0,10] = 0.24
$arr[0,11] = 0.05
$arr[0,12] = 0
$arr[0,13] = 0
$arr[0,14] = 0
$arr[0,15] = 2.93
$arr[0,16] = 24.67
$arr[0,17] = 4.13
$arr[0,18] = 14
$arr[0,19] = 2
$arr[0,20] = 5
$arr[0,21] = 2
$ws.Range("A509:V509").Value = $arr

$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "Entrainement"
$arr[0,1] = 45917
$arr[0,2] = "Global"
$arr[0,3] = "J-3"
$arr[0,4] = "Yoann Martelat"
$arr[0,5] = "center midfield"
$arr[0,6] = "01:36:38"
$arr[0,7] = 6.46
$arr[0,8] = 0.28
$arr[0,9] = 6.17
$arr[0,10] = 0.28
$arr[0,11] = 0
$arr[0,12] = 0
$arr[0,13] = 0
$arr[0,14] = 0
$arr[0,15] = 3.95
$arr[0,16] = 20.42
$arr[0,17] = 4.31
$arr[0,18] = 11
$arr[0,19] = 3
$arr[0,20] = 7
$arr[0,21] = 1
$ws.Range("A510:V510").Value = $arr

$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "Entrainement"
$arr[0,1] = 45917
$arr[0,2] = "Global"
$arr[0,3] = "J-3"
$arr[0,4] = "Emmanuel Valey"
$arr[0,5] = "left forward"
$arr[0,6] = "01:36:28"
$arr[0,7] = 6.55
$arr[0,8] = 0.8
$arr[0,9] = 5.73
$arr[0,10] = 0.55
$arr[0,11] = 0.2
$arr[0,12] = 0.06
$arr[0,13] = 0.01
$arr[0,14] = 5
$arr[0,15] = 3.72
$arr[0,16] = 31.62
$arr[0,17] = 6.28
$arr[0,18] = 71
$arr[0,19] = 24
$arr[0,20] = 46
$arr[0,21] = 17
$ws.Range("A511:V511").Value = $arr

$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "Entrainement"
$arr[0,1] = 45917
$arr[0,2] = "Global"
$arr[0,3] = "J-3"
$arr[0,4] = "Hedi Nasri"
$arr[0,5] = "right back"
$arr[0,6] = "01:37:14"
$arr[0,7] = 6.09
$arr[0,8] = 0.54
$arr[0,9] = 5.54
$arr[0,10] = 0.41
$arr[0,11] = 0.13
$arr[0,12] = 0.01
$arr[0,13] = 0
$arr[0,14] = 2
$arr[0,15] = 3.67
$arr[0,16] = 27.06
$arr[0,17] = 5.37
$arr[0,18] = 28
$arr[0,19] = 15
$arr[0,20] = 24
$arr[0,21] = 2
$ws.Range("A512:V512").Value = $arr

$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "Entrainement"
$arr[0,1] = 45917
$arr[0,2] = "Global"
$arr[0,3] = "J-3"
$arr[0,4] = "Malik Boussaid"
$arr[0,5] = "right back"
$arr[0,6] = "01:39:47"
$arr[0,7] = 6.31
$arr[0,8] = 0.68
$arr[0,9] = 5.62
$arr[0,10] = 0.47
$arr[0,11] = 0.15
$arr[0,12] = 0.06
$arr[0,13] = 0
$arr[0,14] = 5
$arr[0,15] = 3.67
$arr[0,16] = 28.31
$arr[0,17] = 4.68
$arr[0,18] = 32
$arr[0,19] = 6
$arr[0,20] = 19
$arr[0,21] = 5
$ws.Range("A513:V513").Value = $arr

$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "Entrainement"
$arr[0,1] = 45917
$arr[0,2] = "Global"
$arr[0,3] = "J-3"
$arr[0,4] = "Sofiane Belle"
$arr[0,5] = "left forward"
$arr[0,6] = "01:36:37"
$arr[0,7] = 6.72
$arr[0,8] = 0.99
$arr[0,9] = 5.72
$arr[0,10] = 0.57
$arr[0,11] = 0.34
$arr[0,12] = 0.1
$arr[0,13] = 0
$arr[0,14] = 10
$arr[0,15] = 3.99
$arr[0,16] = 28
$arr[0,17] = 4.5
$arr[0,18] = 22
$arr[0,19] = 7
$arr[0,20] = 21
$arr[0,21] = 5
$ws.Range("A514:V514").Value = $arr

# Update sheet view to reflect the new scroll position / selection
$ws.Application.ActiveWindow.ScrollRow = 487
$ws.Range("D520").Select()
